# CSV_Upload.xlsx — add a second sample product (#222), a Size/Color option row,
# a second shipping location row, and a brand-new "Images" sheet; update selections.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Products sheet — add row 3 for the new sample product "#222"
# ---------------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")

$wsProducts.Range("A3").Value = "#222"
$wsProducts.Range("B3").Value = 22222
$wsProducts.Range("C3").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum."""
$wsProducts.Range("D3").Value = "Sample Blah blah blah"
$wsProducts.Range("E3").Value = "New"
$wsProducts.Range("F3").Value = "Iphone"
$wsProducts.Range("G3").Value = 32000
$wsProducts.Range("H3").Value = "Jewellery & Watches"
$wsProducts.Range("I3").Value = "Style"
$wsProducts.Range("J3").Value = "Sanrio"
$wsProducts.Range("K3").Value = "partner2"
$wsProducts.Range("L3").Value = 12
$wsProducts.Range("M3").Value = "sample2"
$wsProducts.Range("N3").Value = 1
$wsProducts.Range("O3").Value = "2.jpg"
$wsProducts.Range("P3").Value = 13

# ---------------------------------------------------------------------------
# 2) Attributes sheet — insert a Size/Small option row for "#111" and append
#    a Color/Indigo option row for "#222"; add the new "Option Image" column
# ---------------------------------------------------------------------------
$wsAttributes = $wb.Worksheets.Item("Attributes")

# insert a fresh row at 3 (pushes the existing "#222/Color/Orange/333" row to 4)
$wsAttributes.Rows.Item(3).Insert()

$wsAttributes.Range("E1").Value = "Option Image"

$wsAttributes.Range("A3").Value = "#111"
$wsAttributes.Range("B3").Value = "Size"
$wsAttributes.Range("C3").Value = "Small"
$wsAttributes.Range("D3").Value = 222

$wsAttributes.Range("A5").Value = "#222"
$wsAttributes.Range("B5").Value = "Color"
$wsAttributes.Range("C5").Value = "Indigo"
$wsAttributes.Range("D5").Value = 444

# the row insert pushed the "locations" validation from A11:A1048576 to
# A12:A1048576; shuffle rows (all currently blank) below our data so it lands
# back on A6:A1048576 without disturbing the real rows 1-5
$wsAttributes.Rows.Item(6).Delete()
$wsAttributes.Rows.Item(6).Delete()
$wsAttributes.Rows.Item(6).Delete()
$wsAttributes.Rows.Item(6).Delete()
$wsAttributes.Rows.Item(6).Delete()
$wsAttributes.Rows.Item(6).Delete()
$wsAttributes.Rows.Item(7).Insert()
$wsAttributes.Rows.Item(7).Insert()
$wsAttributes.Rows.Item(7).Insert()
$wsAttributes.Rows.Item(7).Insert()
$wsAttributes.Rows.Item(7).Insert()
$wsAttributes.Rows.Item(7).Insert()

# ---------------------------------------------------------------------------
# 3) Shipment sheet — append Mindanao/NCR shipping rows for product "#222"
# ---------------------------------------------------------------------------
$wsShipment = $wb.Worksheets.Item("Shipment")

# shuffle currently-blank rows so the "locations" validation (A7:A1048576
# B2:B1048576) becomes A5:A1048576 B2:B1048576, freeing rows 3-4 for data
$wsShipment.Rows.Item(3).Delete()
$wsShipment.Rows.Item(3).Delete()
$wsShipment.Rows.Item(6).Insert()
$wsShipment.Rows.Item(6).Insert()

$wsShipment.Range("A3").Value = "#222"
$wsShipment.Range("B3").Value = "Mindanao"
$wsShipment.Range("C3").Value = 3333

$wsShipment.Range("A4").Value = "#222"
$wsShipment.Range("B4").Value = "NCR"
$wsShipment.Range("C4").Value = 444

# ---------------------------------------------------------------------------
# 4) New "Images" sheet — product number -> image filename lookup
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsImages = $wb.Worksheets.Add($null, $lastSheet)
$wsImages.Name = "Images"

$wsImages.Columns.Item(1).ColumnWidth = 15.33
$wsImages.Columns.Item(2).ColumnWidth = 15.33

$wsImages.Range("A1").Value = "Product Number"
$wsImages.Range("B1").Value = "Product Image File"

$wsImages.Range("A2").Value = "#111"
$wsImages.Range("B2").Value = "1.png"

$wsImages.Range("A3").Value = "#111"
$wsImages.Range("B3").Value = "2.jpg"

$wsImages.Range("A4").Value = "#222"
$wsImages.Range("B4").Value = "1.png"

$wsImages.Range("A5").Value = "#222"
$wsImages.Range("B5").Value = "2.jpg"

$wsImages.Range("A6").Value = "#222"
$wsImages.Range("B6").Value = "2.jpg"

$wsImages.Range("A7:A1048576").Validation.Add(3, 1, 1, "locations")

# ---------------------------------------------------------------------------
# 5) Selections — restore the per-sheet cursor positions and make "Shipment"
#    the active tab (everything else must be selected before the final
#    Activate so only Shipment keeps tabSelected)
# ---------------------------------------------------------------------------
$wsProducts.Range("M20").Select()
$wsAttributes.Range("C11").Select()
$wsImages.Range("B5").Select()

$wsShipment.Activate()
$wsShipment.Range("D7").Select()
